$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header text updates ("Volume ... Number 22" / new report week) ----
$ws.Range("A8").Value = "Volume 30   Number  22"
$ws.Range("C9").Value = "Report Covering the Week  5/29/2023  Through  6/4/2023"

# ---- Plain numeric updates (value only; cell stays numeric, style untouched) ----
$ws.Range("N14").Value = -77.777777777777
$ws.Range("I15").Value = 6
$ws.Range("K15").Value = 50
$ws.Range("L15").Value = 100
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = -60
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -71.428571428571
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -37.5
$ws.Range("I16").Value = 50
$ws.Range("J16").Value = 68
$ws.Range("K16").Value = -26.470588235294
$ws.Range("L16").Value = -24.242424242424
$ws.Range("M16").Value = -43.820224719101
$ws.Range("N16").Value = -86.149584487534
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -57.142857142857
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 18
$ws.Range("H17").Value = -16.666666666666
$ws.Range("I17").Value = 70
$ws.Range("J17").Value = 73
$ws.Range("K17").Value = -4.109589041095
$ws.Range("L17").Value = 4.477611940298
$ws.Range("M17").Value = 48.936170212766
$ws.Range("N17").Value = -61.325966850828
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 27.272727272727
$ws.Range("I18").Value = 70
$ws.Range("J18").Value = 90
$ws.Range("K18").Value = -22.222222222222
$ws.Range("L18").Value = 48.936170212766
$ws.Range("M18").Value = 27.272727272727
$ws.Range("N18").Value = -84.444444444444
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 28.571428571428
$ws.Range("F19").Value = 42
$ws.Range("G19").Value = 27
$ws.Range("H19").Value = 55.555555555555
$ws.Range("I19").Value = 218
$ws.Range("J19").Value = 202
$ws.Range("K19").Value = 7.920792079207
$ws.Range("L19").Value = 28.994082840236
$ws.Range("M19").Value = -1.801801801801
$ws.Range("N19").Value = -48.946135831381
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = 37.5
$ws.Range("I20").Value = 46
$ws.Range("J20").Value = 27
$ws.Range("K20").Value = 70.37037037037
$ws.Range("L20").Value = 48.387096774193
$ws.Range("M20").Value = 170.588235294118
$ws.Range("N20").Value = -88.83495145631
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = -23.076923076923
$ws.Range("F21").Value = 93
$ws.Range("H21").Value = 16.25
$ws.Range("I21").Value = 462
$ws.Range("J21").Value = 465
$ws.Range("K21").Value = -0.645161290322
$ws.Range("L21").Value = 20.626631853785
$ws.Range("M21").Value = 6.451612903225
$ws.Range("N21").Value = -75.094339622641
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("I22").Value = 8
$ws.Range("J22").Value = 13
$ws.Range("K22").Value = -38.461538461538
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -42.857142857142
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = -80
$ws.Range("F23").Value = 6
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = -33.333333333333
$ws.Range("I23").Value = 44
$ws.Range("J23").Value = 56
$ws.Range("K23").Value = -21.428571428571
$ws.Range("L23").Value = -20
$ws.Range("M23").Value = 22.222222222222
$ws.Range("C24").Value = 39
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = 44.444444444444
$ws.Range("F24").Value = 156
$ws.Range("G24").Value = 118
$ws.Range("H24").Value = 32.203389830508
$ws.Range("I24").Value = 770
$ws.Range("J24").Value = 756
$ws.Range("K24").Value = 1.851851851851
$ws.Range("L24").Value = 78.654292343387
$ws.Range("M24").Value = 93.467336683417
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 60
$ws.Range("F25").Value = 31
$ws.Range("G25").Value = 20
$ws.Range("H25").Value = 55
$ws.Range("I25").Value = 128
$ws.Range("J25").Value = 121
$ws.Range("K25").Value = 5.785123966942
$ws.Range("L25").Value = 8.474576271186
$ws.Range("M25").Value = 0
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 100
$ws.Range("I26").Value = 8
$ws.Range("K26").Value = 33.333333333333
$ws.Range("L26").Value = 60
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("I27").Value = 27
$ws.Range("J27").Value = 22
$ws.Range("K27").Value = 22.727272727272
$ws.Range("L27").Value = 28.571428571428
$ws.Range("F28").Value = 2
$ws.Range("I28").Value = 3
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 200
$ws.Range("N28").Value = -85
$ws.Range("F29").Value = 2
$ws.Range("I29").Value = 3
$ws.Range("K29").Value = 50
$ws.Range("L29").Value = 50
$ws.Range("M29").Value = 200
$ws.Range("N29").Value = -85
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 0
$ws.Range("L30").Value = 12.5

# ---- Cells that flip between numeric and text (shared "0" / "***.*" placeholders) ----
# Helper cell (off the used range) staged as Text so its Value pastes as a literal string,
# not an auto-coerced number. Style is fixed up afterwards from a stable donor cell.
$helperText = $ws.Range("Z1")
$helperText.NumberFormat = "@"
$styleDonorText = $ws.Range("C14")   # stable "general/text" style (s=14) donor, never itself retyped
$styleDonorNum  = $ws.Range("I14")   # stable "integer" style (s=15) donor, never itself retyped

$helperText.Value = "0"
$helperText.Copy()
$ws.Range("G14").PasteSpecial(-4163)   # xlPasteValues -> literal text "0"
$styleDonorText.Copy()
$ws.Range("G14").PasteSpecial(-4122)   # xlPasteFormats -> restore numeric-row style 14

$helperText.Value = "***.*"
$helperText.Copy()
$ws.Range("H14").PasteSpecial(-4163)   # xlPasteValues -> literal text "***.*"
$styleDonorText.Copy()
$ws.Range("H14").PasteSpecial(-4122)   # xlPasteFormats -> restore numeric-row style 14

$helperText.Value = "0"
$helperText.Copy()
$ws.Range("G15").PasteSpecial(-4163)   # xlPasteValues -> literal text "0"
$styleDonorText.Copy()
$ws.Range("G15").PasteSpecial(-4122)   # xlPasteFormats -> restore numeric-row style 14

$helperText.Value = "***.*"
$helperText.Copy()
$ws.Range("H15").PasteSpecial(-4163)   # xlPasteValues -> literal text "***.*"
$styleDonorText.Copy()
$ws.Range("H15").PasteSpecial(-4122)   # xlPasteFormats -> restore numeric-row style 14

$helperText.Value = "0"
$helperText.Copy()
$ws.Range("D22").PasteSpecial(-4163)   # xlPasteValues -> literal text "0"
$styleDonorText.Copy()
$ws.Range("D22").PasteSpecial(-4122)   # xlPasteFormats -> restore numeric-row style 14

$helperText.Value = "***.*"
$helperText.Copy()
$ws.Range("E22").PasteSpecial(-4163)   # xlPasteValues -> literal text "***.*"
$styleDonorText.Copy()
$ws.Range("E22").PasteSpecial(-4122)   # xlPasteFormats -> restore numeric-row style 14

$helperText.Value = "0"
$helperText.Copy()
$ws.Range("C26").PasteSpecial(-4163)   # xlPasteValues -> literal text "0"
$styleDonorText.Copy()
$ws.Range("C26").PasteSpecial(-4122)   # xlPasteFormats -> restore numeric-row style 14

$helperText.Value = "0"
$helperText.Copy()
$ws.Range("D26").PasteSpecial(-4163)   # xlPasteValues -> literal text "0"
$styleDonorText.Copy()
$ws.Range("D26").PasteSpecial(-4122)   # xlPasteFormats -> restore numeric-row style 14

$helperText.Value = "***.*"
$helperText.Copy()
$ws.Range("E26").PasteSpecial(-4163)   # xlPasteValues -> literal text "***.*"
$styleDonorText.Copy()
$ws.Range("E26").PasteSpecial(-4122)   # xlPasteFormats -> restore numeric-row style 14

$helperText.Value = "0"
$helperText.Copy()
$ws.Range("G28").PasteSpecial(-4163)   # xlPasteValues -> literal text "0"
$styleDonorText.Copy()
$ws.Range("G28").PasteSpecial(-4122)   # xlPasteFormats -> restore numeric-row style 14

$helperText.Value = "***.*"
$helperText.Copy()
$ws.Range("H28").PasteSpecial(-4163)   # xlPasteValues -> literal text "***.*"
$styleDonorText.Copy()
$ws.Range("H28").PasteSpecial(-4122)   # xlPasteFormats -> restore numeric-row style 14

$helperText.Value = "0"
$helperText.Copy()
$ws.Range("G29").PasteSpecial(-4163)   # xlPasteValues -> literal text "0"
$styleDonorText.Copy()
$ws.Range("G29").PasteSpecial(-4122)   # xlPasteFormats -> restore numeric-row style 14

$helperText.Value = "***.*"
$helperText.Copy()
$ws.Range("H29").PasteSpecial(-4163)   # xlPasteValues -> literal text "***.*"
$styleDonorText.Copy()
$ws.Range("H29").PasteSpecial(-4122)   # xlPasteFormats -> restore numeric-row style 14

$ws.Range("F15").Value = 1
$styleDonorNum.Copy()
$ws.Range("F15").PasteSpecial(-4122)   # xlPasteFormats -> restore style 15

$helperText.Clear()
$excel.CutCopyMode = $false